$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.718.05'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '2.499.97'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D5").Value = "'323.00"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = "'108.55"
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("D7").Value = "'0.523"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.557"
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("D10").Value = "'40.44"
$ws.Range("E10").Value = '  +4.82%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = "'0.125"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = "'18.95"
$ws.Range("E13").Value = '  +3.46%  '
$ws.Range("D14").Value = "'7.18"
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '2.889.62'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '2.493.42'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '47.631.07'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("D20").Value = "'6.60"
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0940'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").Value = "'2.75"
$ws.Range("E22").Value = '  +10.65%  '
$ws.Range("D23").Value = "'70.75"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").Value = "'247.45"
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = "'25.82"
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = '  -4.32%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = "'0.142"
$ws.Range("E30").Value = '  +2.73%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = "'34.73"
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = "'49.80"
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Value = "'19.97"
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").Value = "'5.38"
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = "'1.01"
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'0.0787"
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = "'1.97"
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("D38").Value = "'4.70"
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").Value = "'2.97"
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = "'22.10"
$ws.Range("E41").Value = '  +2.92%  '
$ws.Range("E42").Value = '  -1.65%  '
$ws.Range("D43").Value = "'119.32"
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '2.002.85'
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D46").Value = "'3.07"
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").Value = "'9.04"
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").Value = "'5.13"
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("D51").Value = "'56.35"
$ws.Range("E51").Value = '  +2.12%  '
